$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns for output current calculations
$ws.Range("M1").Value = "Io (NL)"
$ws.Range("N1").Value = "Io at 1k (mA)"
$ws.Range("O1").Value = "Io at 470 (mA)"

# Row 4 (Physical)
$ws.Range("M4").Value = 0
$ws.Range("N4").Formula = "=1000*K4/1000"
$ws.Range("O4").Formula = "=1000*L4/470"

# Row 5 (Sim)
$ws.Range("N5").Formula = "=1000*K5/1000"
$ws.Range("O5").Formula = "=1000*L5/470"

# Column widths (best fit) for new columns
$ws.Columns.Item(14).ColumnWidth = 11
$ws.Columns.Item(15).ColumnWidth = 12

# Update selection to match diff
$ws.Range("O5").Select()
